# Atualização automática de preços de eletricidade
# Updates row 2 of the Spot_PT sheet with the latest daily spot price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spot_PT")

# Day (Excel serial date) -> 2026-02-08
$ws.Range("A2").Value = 46061

# Hourly prices 0h-1h ... 23h-24h
$ws.Range("B2").Value  = 0
$ws.Range("C2").Value  = 0
$ws.Range("D2").Value  = -0.05
$ws.Range("E2").Value  = -0.1
$ws.Range("F2").Value  = -0.1
$ws.Range("G2").Value  = -0.08
$ws.Range("H2").Value  = -0.01
$ws.Range("I2").Value  = 0
$ws.Range("J2").Value  = 0
$ws.Range("K2").Value  = -0.16
$ws.Range("L2").Value  = -0.15
$ws.Range("M2").Value  = -0.15
$ws.Range("N2").Value  = -0.14
$ws.Range("O2").Value  = -0.12
$ws.Range("P2").Value  = -0.15
$ws.Range("Q2").Value  = -0.15
$ws.Range("R2").Value  = -0.14
$ws.Range("S2").Value  = -0.04
$ws.Range("T2").Value  = 1.97
$ws.Range("U2").Value  = 4.19
$ws.Range("V2").Value  = 4.19
$ws.Range("W2").Value  = 4.26
$ws.Range("X2").Value  = 4.19
$ws.Range("Y2").Value  = 1.91
$ws.Range("Z2").Value  = 0.8

# Daily average / slot information
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 3.64
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 4.22
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 3.08
$ws.Range("AG2").Value = "0h-17h"
